$d = $word.ActiveDocument

# Update the date heading paragraph
$d.Content.Find.Execute("2024-04-14 Sunday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2024-04-15 Monday", 2)

# Update each division expression cell-by-cell (row/col addressed so the two
# duplicate "18÷5=" cells each get their own distinct replacement).
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text  = "93÷3="
$t.Cell(1, 2).Range.Text  = "34÷2="
$t.Cell(1, 3).Range.Text  = "16÷5="
$t.Cell(1, 4).Range.Text  = "96÷3="
$t.Cell(1, 5).Range.Text  = "73÷6="

$t.Cell(5, 1).Range.Text  = "64÷8="
$t.Cell(5, 2).Range.Text  = "33÷7="
$t.Cell(5, 3).Range.Text  = "38÷9="
$t.Cell(5, 4).Range.Text  = "63÷9="
$t.Cell(5, 5).Range.Text  = "88÷2="

$t.Cell(9, 1).Range.Text  = "29÷5="
$t.Cell(9, 2).Range.Text  = "72÷9="
$t.Cell(9, 3).Range.Text  = "63÷6="
$t.Cell(9, 4).Range.Text  = "91÷3="
$t.Cell(9, 5).Range.Text  = "63÷5="

$t.Cell(13, 1).Range.Text = "13÷4="
$t.Cell(13, 2).Range.Text = "31÷2="
$t.Cell(13, 3).Range.Text = "88÷8="
$t.Cell(13, 4).Range.Text = "76÷3="
$t.Cell(13, 5).Range.Text = "74÷9="

$t.Cell(17, 1).Range.Text = "26÷8="
$t.Cell(17, 2).Range.Text = "65÷5="
$t.Cell(17, 3).Range.Text = "59÷5="
$t.Cell(17, 4).Range.Text = "66÷2="
$t.Cell(17, 5).Range.Text = "51÷8="
